# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. immediately
#    before the existing "2022-Q1" sheet) and populate it with the fund
#    holdings data for that quarter (same layout as the other quarter
#    sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
# 2. Update the "总计" (summary) sheet so it gains a new top row for
#    2022-Q4 and the previously-existing rows shift down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create + place the new "2022-Q4" worksheet
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q1")
$q4 = $wb.Worksheets.Add($anchor)
$q4.Name = "2022-Q4"

# Headers (row 1, columns B:H) -- values first, styling copied afterwards
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Row 2 -- 天弘中证1000指数增强C
$q4.Range("A2").Value = 0
$q4.Range("C2").Value = "天弘中证1000指数增强C"
$q4.Range("H2").Value = 1

# Row 3 -- 天弘中证1000指数增强A
$q4.Range("A3").Value = 1
$q4.Range("C3").Value = "天弘中证1000指数增强A"
$q4.Range("H3").Value = 1

# Columns that hold numeric-looking values which must stay TEXT (fund
# code keeps its leading zero, the decimal figures keep trailing zeros):
# B2:B3 (基金代码), D2:E3 (基金规模/股票总仓位), F2:F3 (仓位占比), G2:G3 (持有市值)
$q4.Range("B2:B3").NumberFormat = "@"
$q4.Range("D2:G3").NumberFormat = "@"

$q4.Range("B2").Value = "014202"
$q4.Range("D2").Value = "6.60"
$q4.Range("E2").Value = "94.11"
$q4.Range("F2").Value = "1.66"
$q4.Range("G2").Value = "0.1096"

$q4.Range("B3").Value = "014201"
$q4.Range("D3").Value = "3.86"
$q4.Range("E3").Value = "94.11"
$q4.Range("F3").Value = "1.66"
$q4.Range("G3").Value = "0.0641"

# Copy the bold/bordered header style from an existing quarter sheet onto
# the header row and the index column (A2:A3) of the new sheet.
$styleSrc = $wb.Worksheets.Item("2021-Q3")
$styleSrc.Range("B1:H1").Copy() | Out-Null
$q4.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$styleSrc.Range("A2:A3").Copy() | Out-Null
$q4.Range("A2:A3").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Step 2: shift the "总计" summary rows down and insert the 2022-Q4 row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.14

$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.14

$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.05

$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.28

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.17

# New row 6 needs the same bold/bordered style as the existing index
# column cells (A2:A5) applied to A6, plus the sequential index value.
$total.Range("A6").Value = 4
$total.Range("A2").Copy() | Out-Null
$total.Range("A6").PasteSpecial(-4122) | Out-Null

Write-Output "2022-Q4 sheet + 总计 summary updated"
